$wb = $excel.ActiveWorkbook

# --- 1. Produk sheet: insert a new column F for 'ongkosKirimBeli' ---
$produk = $wb.Worksheets.Item("Produk")

# Insert a new column before column F (6th column), shifting existing F..P to G..Q
$produk.Columns.Item(6).Insert()

# Header cell for the new column
$produk.Cells.Item(1, 6).Value = "ongkosKirimBeli"

# Approximate the width of the neighboring column E for the new column F
$produk.Columns.Item(6).ColumnWidth = $produk.Columns.Item(5).ColumnWidth

# --- 2. Update active sheet / selection state ---
# Previously sheet "PengeluaranBarang" was the tab-selected / active sheet;
# now "Produk" becomes the active sheet with F2 selected.
$produk.Activate()
$produk.Range("F2").Select()
